$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1. Fix wording in first (Descripcion) paragraph:
#    "entendimiento de dinámica de la tierra" -> "entendimiento dinámico de la tierra"
Replace-Text "entendimiento de dinámica de la tierra" "entendimiento dinámico de la tierra"

# 2. Bibliography clean-up: each reference list currently lives as a single
#    paragraph whose citations are separated only by a lone-space run. Turn
#    every one of those separating spaces into a real paragraph break so
#    each citation becomes its own paragraph (still styled "TextBody", which
#    Word inherits automatically for the newly split-off paragraph).

Replace-Text `
  "Referencias Básicas. Odum, Eugene P., Barret H, Garry W, y Aguilar Orte" `
  "Referencias Básicas.^pOdum, Eugene P., Barret H, Garry W, y Aguilar Orte"

Replace-Text `
  " Ecología. Madrid: CENGAGE Learning Latin America. Smith, Thomas M., y Smith, Robert Leo. (2007). Eco" `
  " Ecología. Madrid: CENGAGE Learning Latin America.^pSmith, Thomas M., y Smith, Robert Leo. (2007). Eco"

Replace-Text `
  " Leo. (2007). Ecología. Madrid: Pearson Educación. Sarmiento, Guillermo (1980). Los Ecosistemas y la " `
  " Leo. (2007). Ecología. Madrid: Pearson Educación.^pSarmiento, Guillermo (1980). Los Ecosistemas y la "

Replace-Text `
  "orial Blume, S.A. Milanesat 21-23 08017 Barcelona. Vázquez Torre Ana María Guadalupe (1993). Ecología" `
  "orial Blume, S.A. Milanesat 21-23 08017 Barcelona.^pVázquez Torre Ana María Guadalupe (1993). Ecología"

Replace-Text `
  "Referencias Complementarias. Boege Ek del Val, Karina. (2012). Ecología y Evolu" `
  "Referencias Complementarias.^pBoege Ek del Val, Karina. (2012). Ecología y Evolu"

Replace-Text `
  "s Bióticas. Argentina: Fondo de Cultura Económica. Escolástico León, Consuelo., Cabildo Miranda, Marí" `
  "s Bióticas. Argentina: Fondo de Cultura Económica.^pEscolástico León, Consuelo., Cabildo Miranda, Marí"

Replace-Text `
  "troducción Organismos y Poblaciones. UNED: Madrid. Escolástico León, Consuelo., Cabildo Miranda, Marí" `
  "troducción Organismos y Poblaciones. UNED: Madrid.^pEscolástico León, Consuelo., Cabildo Miranda, Marí"

Replace-Text `
  "logía II: Comunidades y Ecosistemas. UNED: Madrid. Samo Lumbreras, Antonio José., Garmendia Salvador," `
  "logía II: Comunidades y Ecosistemas. UNED: Madrid.^pSamo Lumbreras, Antonio José., Garmendia Salvador,"

Replace-Text `
  "tica a la Ecología. Pearson Prentice Hall: Madrid. Tarbuck, Edward J., y Lutgens Frederick K. (2005)." `
  "tica a la Ecología. Pearson Prentice Hall: Madrid.^pTarbuck, Edward J., y Lutgens Frederick K. (2005)."

Replace-Text `
  "ncias de la Tierra. Madrid: Pearson Prentice Hall. Yánez, Patricio. (2014). Ecología y Biodiversidad " `
  "ncias de la Tierra. Madrid: Pearson Prentice Hall.^pYánez, Patricio. (2014). Ecología y Biodiversidad "

Replace-Text `
  "d desde el neotrópico: UNIBE/UIDE: Quito, Ecuador. Magurran, A. E. (1988). Ecological diversity and i" `
  "d desde el neotrópico: UNIBE/UIDE: Quito, Ecuador.^pMagurran, A. E. (1988). Ecological diversity and i"

Replace-Text `
  "n University Press, Princeton, New Jersey, 179 pp. Alcolado, P.M. (1984). Conceptos e índices relacio" `
  "n University Press, Princeton, New Jersey, 179 pp.^pAlcolado, P.M. (1984). Conceptos e índices relacio"

Replace-Text `
  "dos con la diversidad. Avicennia, 1998, 8/9: 7-21. Alcolado, P.M. (1984). Utilidad de algunos Índices" `
  "dos con la diversidad. Avicennia, 1998, 8/9: 7-21.^pAlcolado, P.M. (1984). Utilidad de algunos Índices"

Write-Host "edit complete"
